# Update "F" column ("想去人数" / want-to-go count) values on three sheets:
# 展览 (Exhibition), 演出 (Performance), 全部类型 (All types).
# 本地生活 (Local life) is not touched.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 206
$ws1.Cells.Item(6, 6).Value = 1369
$ws1.Cells.Item(7, 6).Value = 78
$ws1.Cells.Item(10, 6).Value = 477
$ws1.Cells.Item(11, 6).Value = 838
$ws1.Cells.Item(12, 6).Value = 542
$ws1.Cells.Item(13, 6).Value = 759
$ws1.Cells.Item(14, 6).Value = 331
$ws1.Cells.Item(15, 6).Value = 502
$ws1.Cells.Item(17, 6).Value = 1071
$ws1.Cells.Item(18, 6).Value = 523
$ws1.Cells.Item(22, 6).Value = 253
$ws1.Cells.Item(26, 6).Value = 474
$ws1.Cells.Item(28, 6).Value = 363

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 390
$ws2.Cells.Item(4, 6).Value = 57
$ws2.Cells.Item(10, 6).Value = 164
$ws2.Cells.Item(11, 6).Value = 162
$ws2.Cells.Item(14, 6).Value = 16

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 206
$ws4.Cells.Item(7, 6).Value = 1369
$ws4.Cells.Item(9, 6).Value = 78
$ws4.Cells.Item(10, 6).Value = 390
$ws4.Cells.Item(12, 6).Value = 57
$ws4.Cells.Item(16, 6).Value = 477
$ws4.Cells.Item(17, 6).Value = 838
$ws4.Cells.Item(18, 6).Value = 542
$ws4.Cells.Item(19, 6).Value = 759
$ws4.Cells.Item(20, 6).Value = 331
$ws4.Cells.Item(21, 6).Value = 502
$ws4.Cells.Item(23, 6).Value = 1071
$ws4.Cells.Item(24, 6).Value = 523
$ws4.Cells.Item(31, 6).Value = 164
$ws4.Cells.Item(32, 6).Value = 253
$ws4.Cells.Item(35, 6).Value = 162
$ws4.Cells.Item(39, 6).Value = 16
$ws4.Cells.Item(40, 6).Value = 474
$ws4.Cells.Item(42, 6).Value = 363
